$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 0.1541501976284585
$ws.Cells.Item(2, 3).Value = 0.5889328063241107
$ws.Cells.Item(2, 10).Value = 0.02371541501976284
$ws.Cells.Item(2, 16).Value = 0.1264822134387352
$ws.Cells.Item(2, 19).Value = 0.1067193675889328
$ws.Cells.Item(3, 2).Value = 0.01257861635220126
$ws.Cells.Item(3, 3).Value = 0.05031446540880503
$ws.Cells.Item(3, 10).Value = 0.01886792452830189
$ws.Cells.Item(3, 16).Value = 0.7672955974842768
$ws.Cells.Item(3, 19).Value = 0.1509433962264151
$ws.Cells.Item(4, 10).Value = 0.1081081081081081
$ws.Cells.Item(4, 16).Value = 0.5135135135135135
$ws.Cells.Item(4, 19).Value = 0.3783783783783784
$ws.Cells.Item(6, 2).Value = 0.04761904761904762
$ws.Cells.Item(6, 4).Value = 0.0119047619047619
$ws.Cells.Item(6, 6).Value = 0.08333333333333333
$ws.Cells.Item(6, 10).Value = 0.2341269841269841
$ws.Cells.Item(6, 15).Value = 0.01587301587301587
$ws.Cells.Item(6, 17).Value = 0.1746031746031746
$ws.Cells.Item(6, 18).Value = 0.07539682539682539
$ws.Cells.Item(6, 19).Value = 0.3571428571428572
$ws.Cells.Item(7, 2).Value = 0.116751269035533
$ws.Cells.Item(7, 4).Value = 0.02538071065989848
$ws.Cells.Item(7, 6).Value = 0.04568527918781726
$ws.Cells.Item(7, 10).Value = 0.116751269035533
$ws.Cells.Item(7, 15).Value = 0.02538071065989848
$ws.Cells.Item(7, 17).Value = 0.1319796954314721
$ws.Cells.Item(7, 18).Value = 0.116751269035533
$ws.Cells.Item(7, 19).Value = 0.4213197969543147
$ws.Cells.Item(8, 2).Value = 0.07520325203252033
$ws.Cells.Item(8, 4).Value = 0.02032520325203252
$ws.Cells.Item(8, 6).Value = 0.06504065040650407
$ws.Cells.Item(8, 10).Value = 0.08333333333333333
$ws.Cells.Item(8, 15).Value = 0.01626016260162602
$ws.Cells.Item(8, 17).Value = 0.1788617886178862
$ws.Cells.Item(8, 18).Value = 0.1036585365853658
$ws.Cells.Item(8, 19).Value = 0.4573170731707317
$ws.Cells.Item(9, 2).Value = 0.07100591715976332
$ws.Cells.Item(9, 4).Value = 0.005917159763313609
$ws.Cells.Item(9, 6).Value = 0.106508875739645
$ws.Cells.Item(9, 10).Value = 0.04733727810650887
$ws.Cells.Item(9, 15).Value = 0.01775147928994083
$ws.Cells.Item(9, 17).Value = 0.1893491124260355
$ws.Cells.Item(9, 18).Value = 0.1242603550295858
$ws.Cells.Item(9, 19).Value = 0.4378698224852071
$ws.Cells.Item(10, 2).Value = 0.1129597197898424
$ws.Cells.Item(10, 4).Value = 0.01663747810858144
$ws.Cells.Item(10, 5).Value = 0.001751313485113835
$ws.Cells.Item(10, 6).Value = 0.08844133099824869
$ws.Cells.Item(10, 10).Value = 0.08143607705779335
$ws.Cells.Item(10, 15).Value = 0.01663747810858144
$ws.Cells.Item(10, 17).Value = 0.1961471103327496
$ws.Cells.Item(10, 18).Value = 0.09194395796847636
$ws.Cells.Item(10, 19).Value = 0.3940455341506129
$ws.Cells.Item(11, 7).Value = 0.1447368421052632
$ws.Cells.Item(11, 10).Value = 0.08552631578947369
$ws.Cells.Item(11, 11).Value = 0.194078947368421
$ws.Cells.Item(11, 12).Value = 0.555921052631579
$ws.Cells.Item(11, 19).Value = 0.01973684210526316
$ws.Cells.Item(12, 7).Value = 0.7634408602150538
$ws.Cells.Item(12, 10).Value = 0.1612903225806452
$ws.Cells.Item(12, 11).Value = 0.005376344086021506
$ws.Cells.Item(12, 12).Value = 0.02150537634408602
$ws.Cells.Item(12, 19).Value = 0.04838709677419355
$ws.Cells.Item(13, 6).Value = 0.01818181818181818
$ws.Cells.Item(13, 7).Value = 0.5454545454545454
$ws.Cells.Item(13, 10).Value = 0.3454545454545455
$ws.Cells.Item(13, 19).Value = 0.09090909090909091
$ws.Cells.Item(14, 10).Value = 0.6666666666666666
$ws.Cells.Item(14, 19).Value = 0.3333333333333333
$ws.Cells.Item(15, 6).Value = 0.025
$ws.Cells.Item(15, 8).Value = 0.1958333333333333
$ws.Cells.Item(15, 9).Value = 0.07916666666666666
$ws.Cells.Item(15, 10).Value = 0.3541666666666667
$ws.Cells.Item(15, 11).Value = 0.06666666666666667
$ws.Cells.Item(15, 13).Value = 0.01666666666666667
$ws.Cells.Item(15, 15).Value = 0.05
$ws.Cells.Item(15, 19).Value = 0.2125
$ws.Cells.Item(16, 6).Value = 0.01796407185628742
$ws.Cells.Item(16, 8).Value = 0.2035928143712575
$ws.Cells.Item(16, 9).Value = 0.05389221556886228
$ws.Cells.Item(16, 10).Value = 0.437125748502994
$ws.Cells.Item(16, 11).Value = 0.1317365269461078
$ws.Cells.Item(16, 13).Value = 0.01197604790419162
$ws.Cells.Item(16, 15).Value = 0.04790419161676647
$ws.Cells.Item(16, 19).Value = 0.09580838323353294
$ws.Cells.Item(17, 6).Value = 0.01932367149758454
$ws.Cells.Item(17, 8).Value = 0.2077294685990338
$ws.Cells.Item(17, 9).Value = 0.06763285024154589
$ws.Cells.Item(17, 10).Value = 0.4323671497584541
$ws.Cells.Item(17, 11).Value = 0.0748792270531401
$ws.Cells.Item(17, 13).Value = 0.01449275362318841
$ws.Cells.Item(17, 14).Value = 0.002415458937198068
$ws.Cells.Item(17, 15).Value = 0.07246376811594203
$ws.Cells.Item(17, 19).Value = 0.108695652173913
$ws.Cells.Item(18, 6).Value = 0.01382488479262673
$ws.Cells.Item(18, 8).Value = 0.2258064516129032
$ws.Cells.Item(18, 9).Value = 0.08755760368663594
$ws.Cells.Item(18, 10).Value = 0.3640552995391705
$ws.Cells.Item(18, 11).Value = 0.1105990783410138
$ws.Cells.Item(18, 13).Value = 0.02304147465437788
$ws.Cells.Item(18, 15).Value = 0.08755760368663594
$ws.Cells.Item(18, 19).Value = 0.08755760368663594
$ws.Cells.Item(19, 6).Value = 0.02312599681020734
$ws.Cells.Item(19, 8).Value = 0.2129186602870813
$ws.Cells.Item(19, 9).Value = 0.07177033492822966
$ws.Cells.Item(19, 10).Value = 0.3492822966507177
$ws.Cells.Item(19, 11).Value = 0.1180223285486443
$ws.Cells.Item(19, 13).Value = 0.02791068580542265
$ws.Cells.Item(19, 14).Value = 0.001594896331738437
$ws.Cells.Item(19, 15).Value = 0.08373205741626795
$ws.Cells.Item(19, 19).Value = 0.1116427432216906
